$d = $word.ActiveDocument

$d.Content.Find.Execute("Start time: 2017-12-27 18:45:29", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Start time: 2018-01-31 12:39:00", 2)

$d.Content.Find.Execute("End time: 2017-12-27 18:53:41", $true, $false, $false, $false, $false,
                         $true, 1, $false, "End time: 2018-01-31 12:45:53", 2)

$d.Content.Find.Execute("Duration: 8.19 mins", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Duration: 6.87 mins", 2)
